# Auto-generated edit script applying the cryptos.xlsx data refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.638.35'
$ws.Range("E2").Value = '  +0.78%  '

$ws.Range("D3").Value = '2.614.77'
$ws.Range("E3").Value = '  +0.50%  '

$ws.Range("E4").Value = '  -0.41%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '513.24'
$ws.Range("E5").Value = '  +1.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.61'
$ws.Range("E6").Value = '  -1.27%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").Value = '2.629.44'
$ws.Range("E9").Value = '  -0.37%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.81'
$ws.Range("E10").Value = '  +4.85%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.104'
$ws.Range("E11").Value = '  -0.30%  '

$ws.Range("E12").Value = '  +1.20%  '

$ws.Range("E13").Value = '  +1.94%  '

$ws.Range("D14").Value = '3.077.38'
$ws.Range("E14").Value = '  -0.20%  '

$ws.Range("D15").Value = '60.650.98'
$ws.Range("E15").Value = '  +0.57%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.67'
$ws.Range("E16").Value = '  +0.08%  '

$ws.Range("E17").Value = '  +0.63%  '

$ws.Range("D18").Value = '2.634.14'
$ws.Range("E18").Value = '  -0.23%  '

$ws.Range("E19").Value = '  -0.81%  '

$ws.Range("E20").Value = '  +3.69%  '

$ws.Range("E21").Value = '  +1.45%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.18'
$ws.Range("E22").Value = '  +0.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.67'
$ws.Range("E24").Value = '  +1.08%  '

$ws.Range("E25").Value = '  +0.11%  '

$ws.Range("D26").Value = '2.736.32'
$ws.Range("E26").Value = '  -0.95%  '

$ws.Range("E27").Value = '  +0.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +0.23%  '

$ws.Range("D29").Value = '0.0₃0847'
$ws.Range("E29").Value = '  -1.39%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.36'
$ws.Range("E30").Value = '  -2.38%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.46'
$ws.Range("E32").Value = '  -0.12%  '

$ws.Range("B33").Value = 'Aptos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.96'
$ws.Range("E33").Value = '  +3.82%  '

$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.58'
$ws.Range("E34").Value = '  +0.81%  '

$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '152.04'
$ws.Range("E35").Value = '  -2.54%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.00'
$ws.Range("E36").Value = '  -0.66%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.19'
$ws.Range("E37").Value = '  -1.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.880'
$ws.Range("E38").Value = '  +5.42%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.49'
$ws.Range("E39").Value = '  +0.00%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.850'
$ws.Range("E40").Value = '  +0.34%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.31'
$ws.Range("E41").Value = '  +2.68%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.75'
$ws.Range("E42").Value = '  -0.20%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '292.04'
$ws.Range("E43").Value = '  -5.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.101'
$ws.Range("E44").Value = '  +1.31%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.625'
$ws.Range("E45").Value = '  -0.94%  '

$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.998'
$ws.Range("E46").Value = '  +0.52%  '

$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0556'
$ws.Range("E47").Value = '  -2.67%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.88'
$ws.Range("E48").Value = '  -0.38%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.93'
$ws.Range("E49").Value = '  +2.03%  '

$ws.Range("E50").Value = '  -0.27%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.30'
$ws.Range("E51").Value = '  +0.26%  '
